# Add the new "MPI" worksheet after the existing "OpenMP" sheet and
# populate it with the weak/strong-scaling SVD/DGEMM/HDF5 comparison data.

$wb = $excel.ActiveWorkbook

$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "MPI"

# --- Header row (row 1) ------------------------------------------------
# Written in this specific order so the shared-string table is built up
# with the same index assignment as the source workbook.
$ws.Range("A1").Value = "nodes"
$ws.Range("B1").Value = "cores/node"
$ws.Range("C1").Value = "processes"
$ws.Range("D1").Value = "matrix dim."
$ws.Range("J1").Value = "singularity DGEMM (s)"
$ws.Range("I1").Value = "singularity SVD (s)"
$ws.Range("E1").Value = "native HDF5 (s)"
$ws.Range("F1").Value = "native SVD (s)"
$ws.Range("G1").Value = "native DGEMM (s)"
$ws.Range("H1").Value = "singularity HDF5 (s)"
$ws.Range("K1").Value = "HDF5 overhead"
$ws.Range("L1").Value = "SVD overhead"
$ws.Range("M1").Value = "DGEMM overhead"

# --- Data row (row 2) ---------------------------------------------------
$ws.Range("A2").Value = 2
$ws.Range("B2").Value = 8
$ws.Range("C2").Formula = "=A2*B2"
$ws.Range("D2").Value = 5000

$ws.Range("E2").Value = 1.8913409999999999
$ws.Range("F2").Value = 68.060210999999995
$ws.Range("G2").Value = 2.5266380000000002
$ws.Range("H2").Value = 2.4477449999999998
$ws.Range("I2").Value = 99.785908000000006
$ws.Range("J2").Value = 12.625704000000001

$ws.Range("K2").Formula = "=H2/E2"
$ws.Range("L2").Formula = "=I2/F2"
$ws.Range("M2").Formula = "=J2/G2"

# Numeric columns E2:M2 use the "0.00" number format.
$ws.Range("E2:M2").NumberFormat = "0.00"

# --- Column widths (best-fit approximation) -----------------------------
$widths = @(5.88671875, 10.33203125, 9, 10.21875, 13.21875, 12, 15.6640625, 16.44140625, 15.33203125, 19, 13.5546875, 12.33203125, 16)
for ($i = 0; $i -lt $widths.Length; $i++) {
    $ws.Columns.Item($i + 1).ColumnWidth = $widths[$i] - 5/6
}

# --- Page setup (match the other sheet's print settings) ---------------
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

# --- Selection / active sheet -------------------------------------------
[void]$ws.Range("E2:M2").Select()
